# "Generate Report for Handback" - refresh the localization-status report
# after a successful handback: update status text, clear the stale-handback
# warning, bump the handback timestamps, and widen a couple of report
# columns so the (now longer) status text / (now shorter) error text fit.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This string is shared by Overview!E2 (zh-cn column), Overview!F2
#    (de-de column), zh-cn!C2 and de-de!C2 (the Status columns).
# ---------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$ws_overview.Range("E2").Value = $newStatus
$ws_overview.Range("F2").Value = $newStatus
$ws_zhcn.Range("C2").Value = $newStatus
$ws_dede.Range("C2").Value = $newStatus

# ---------------------------------------------------------------------
# 2) Latest Handback DateTime bumped to the new handback run's timestamp.
# ---------------------------------------------------------------------
$ws_zhcn.Range("K2").Value = "2016-08-12 23:03:16"
$ws_dede.Range("K2").Value = "2016-08-12 23:03:25"

# ---------------------------------------------------------------------
# 3) Error Detail cleared now that the handback file is up to date.
# ---------------------------------------------------------------------
$ws_zhcn.Range("P2").Value = ""
$ws_dede.Range("P2").Value = ""

# ---------------------------------------------------------------------
# 4) Widen the "Status" column (Overview E:F, zh-cn/de-de C) to fit the
#    longer status text, and narrow the "Error Detail" column (P) now
#    that it no longer needs to hold the long warning message.
# ---------------------------------------------------------------------
$ws_overview.Columns.Item(5).ColumnWidth = 29.1666666666667
$ws_overview.Columns.Item(6).ColumnWidth = 29.1666666666667

$ws_zhcn.Columns.Item(3).ColumnWidth = 29.1666666666667
$ws_zhcn.Columns.Item(16).ColumnWidth = 12.8333333333333

$ws_dede.Columns.Item(3).ColumnWidth = 29.1666666666667
$ws_dede.Columns.Item(16).ColumnWidth = 12.8333333333333
